# Scheduled runner refresh: update cached Universalis market-board pricing
# columns (currentAveragePrice* / LevePrice* / LeveProfit*) across the per-job
# profit sheets, matching the latest market data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 335.04544
$ws.Range("I33").Value = 253.8421
$ws.Range("J33").Value = 849.3333
$ws.Range("K33").Value = 253.8421
$ws.Range("L33").Value = 849.3333
$ws.Range("M33").Value = -24.84209999999999
$ws.Range("N33").Value = -1307.3333
$ws.Range("H112").Value = 1708.7142
$ws.Range("J112").Value = 2036
$ws.Range("L112").Value = 6108
$ws.Range("N112").Value = -8324
$ws.Range("H137").Value = 517771.16
$ws.Range("I137").Value = 767.4706
$ws.Range("J137").Value = 1050441.6
$ws.Range("K137").Value = 2302.4118
$ws.Range("L137").Value = 3151324.8
$ws.Range("M137").Value = 247.5882000000001
$ws.Range("N137").Value = -3156424.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1719.2051
$ws.Range("I2").Value = 1686.2667
$ws.Range("K2").Value = 1686.2667
$ws.Range("M2").Value = -1573.2667
$ws.Range("H32").Value = 5752590.5
$ws.Range("I32").Value = 6101468
$ws.Range("J32").Value = 30999.4
$ws.Range("K32").Value = 6101468
$ws.Range("L32").Value = 30999.4
$ws.Range("M32").Value = -6101181
$ws.Range("N32").Value = -31573.4
$ws.Range("H37").Value = 43157.58
$ws.Range("I37").Value = 8999.375
$ws.Range("K37").Value = 8999.375
$ws.Range("M37").Value = -8726.375
$ws.Range("H45").Value = 8678.666999999999
$ws.Range("I45").Value = 6029.857
$ws.Range("K45").Value = 6029.857
$ws.Range("M45").Value = -5652.857
$ws.Range("H74").Value = 2503.1428
$ws.Range("I74").Value = 2075.1191
$ws.Range("J74").Value = 3359.1904
$ws.Range("K74").Value = 2075.1191
$ws.Range("L74").Value = 3359.1904
$ws.Range("M74").Value = -1201.1191
$ws.Range("N74").Value = -5107.190399999999
$ws.Range("H77").Value = 2503.1428
$ws.Range("I77").Value = 2075.1191
$ws.Range("J77").Value = 3359.1904
$ws.Range("K77").Value = 10375.5955
$ws.Range("L77").Value = 16795.952
$ws.Range("M77").Value = -6007.595499999999
$ws.Range("N77").Value = -25531.952
$ws.Range("H102").Value = 3515.7273
$ws.Range("I102").Value = 2775.5557
$ws.Range("K102").Value = 2775.5557
$ws.Range("M102").Value = -1153.5557
$ws.Range("H110").Value = 659.8
$ws.Range("I110").Value = 660.087
$ws.Range("K110").Value = 660.087
$ws.Range("M110").Value = 1384.913
$ws.Range("H116").Value = 1719.2051
$ws.Range("I116").Value = 1686.2667
$ws.Range("K116").Value = 1686.2667
$ws.Range("M116").Value = 607.7333000000001
$ws.Range("H132").Value = 3137.4285
$ws.Range("I132").Value = 2299.3333
$ws.Range("J132").Value = 5651.7144
$ws.Range("K132").Value = 6897.999899999999
$ws.Range("L132").Value = 16955.1432
$ws.Range("M132").Value = -4367.999899999999
$ws.Range("N132").Value = -22015.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1719.2051
$ws.Range("I3").Value = 1686.2667
$ws.Range("K3").Value = 1686.2667
$ws.Range("M3").Value = -1572.2667
$ws.Range("H20").Value = 3171
$ws.Range("I20").Value = 2646.7693
$ws.Range("J20").Value = 4874.75
$ws.Range("K20").Value = 2646.7693
$ws.Range("L20").Value = 4874.75
$ws.Range("M20").Value = -2399.7693
$ws.Range("N20").Value = -5368.75
$ws.Range("H105").Value = 2862.9268
$ws.Range("I105").Value = 2568.5
$ws.Range("K105").Value = 2568.5
$ws.Range("M105").Value = -821.5
$ws.Range("H128").Value = 16399.8
$ws.Range("I128").Value = 16399.8
$ws.Range("K128").Value = 49199.39999999999
$ws.Range("M128").Value = -46709.39999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1998.8572
$ws.Range("I58").Value = 1530.3158
$ws.Range("K58").Value = 1530.3158
$ws.Range("M58").Value = -1327.3158
$ws.Range("H99").Value = 8216.777
$ws.Range("J99").Value = 7900
$ws.Range("L99").Value = 7900
$ws.Range("N99").Value = -10896
$ws.Range("H126").Value = 8216.777
$ws.Range("J126").Value = 7900
$ws.Range("L126").Value = 23700
$ws.Range("N126").Value = -28640
$ws.Range("H132").Value = 6758652.5
$ws.Range("I132").Value = 1771.091
$ws.Range("K132").Value = 5313.272999999999
$ws.Range("M132").Value = -2783.272999999999
$ws.Range("H136").Value = 1998.8572
$ws.Range("I136").Value = 1530.3158
$ws.Range("K136").Value = 4590.9474
$ws.Range("M136").Value = -2040.9474
$ws.Range("H138").Value = 40813
$ws.Range("I138").Value = 48209
$ws.Range("K138").Value = 48209
$ws.Range("M138").Value = -43069
$ws.Range("H139").Value = 83000
$ws.Range("J139").Value = 83000
$ws.Range("L139").Value = 83000
$ws.Range("N139").Value = -93280
$ws.Range("H140").Value = 92439.96000000001
$ws.Range("J140").Value = 92439.96000000001
$ws.Range("L140").Value = 92439.96000000001
$ws.Range("N140").Value = -102799.96
$ws.Range("H141").Value = 217962
$ws.Range("J141").Value = 217962
$ws.Range("L141").Value = 217962
$ws.Range("N141").Value = -228322

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3153.375
$ws.Range("I5").Value = 746.7143
$ws.Range("K5").Value = 2240.1429
$ws.Range("M5").Value = -2128.1429
$ws.Range("H38").Value = 30.61111
$ws.Range("I38").Value = 23.3
$ws.Range("J38").Value = 39.75
$ws.Range("K38").Value = 69.90000000000001
$ws.Range("L38").Value = 119.25
$ws.Range("M38").Value = 277.1
$ws.Range("N38").Value = -813.25
$ws.Range("H56").Value = 8838.333000000001
$ws.Range("I56").Value = 8838.333000000001
$ws.Range("K56").Value = 8838.333000000001
$ws.Range("M56").Value = -8308.333000000001
$ws.Range("H113").Value = 1375.6428
$ws.Range("I113").Value = 949.3333
$ws.Range("J113").Value = 1577.579
$ws.Range("K113").Value = 2847.9999
$ws.Range("L113").Value = 4732.737
$ws.Range("M113").Value = -677.9998999999998
$ws.Range("N113").Value = -9072.737000000001
$ws.Range("H129").Value = 1663.0476
$ws.Range("I129").Value = 1094.5834
$ws.Range("J129").Value = 2421
$ws.Range("K129").Value = 3283.7502
$ws.Range("L129").Value = 7263
$ws.Range("M129").Value = 1716.2498
$ws.Range("N129").Value = -17263
$ws.Range("H131").Value = 2294.147
$ws.Range("I131").Value = 2174
$ws.Range("K131").Value = 6522
$ws.Range("M131").Value = -1482
$ws.Range("H135").Value = 3153.375
$ws.Range("I135").Value = 746.7143
$ws.Range("K135").Value = 6720.428699999999
$ws.Range("M135").Value = -4185.428699999999
$ws.Range("H137").Value = 2636.6667
$ws.Range("J137").Value = 2000.6666
$ws.Range("L137").Value = 6001.9998
$ws.Range("N137").Value = -16201.9998
$ws.Range("H141").Value = 4999.5
$ws.Range("I141").Value = 4999.5
$ws.Range("K141").Value = 14998.5
$ws.Range("M141").Value = -9818.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3817.1428
$ws.Range("I102").Value = 2824.9167
$ws.Range("K102").Value = 2824.9167
$ws.Range("M102").Value = -1202.9167
$ws.Range("H117").Value = 59436
$ws.Range("J117").Value = 59436
$ws.Range("L117").Value = 59436
$ws.Range("N117").Value = -66320
$ws.Range("H132").Value = 209643.86
$ws.Range("I132").Value = 310228.25
$ws.Range("J132").Value = 3181.158
$ws.Range("K132").Value = 930684.75
$ws.Range("L132").Value = 9543.474
$ws.Range("M132").Value = -928154.75
$ws.Range("N132").Value = -14603.474

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4675.6665
$ws.Range("I122").Value = 4280.6875
$ws.Range("J122").Value = 5250.1816
$ws.Range("K122").Value = 12842.0625
$ws.Range("L122").Value = 15750.5448
$ws.Range("M122").Value = -10392.0625
$ws.Range("N122").Value = -20650.5448
$ws.Range("H132").Value = 526064.4399999999
$ws.Range("I132").Value = 707457.1
$ws.Range("J132").Value = 3226.8235
$ws.Range("K132").Value = 2122371.3
$ws.Range("L132").Value = 9680.470499999999
$ws.Range("M132").Value = -2119841.3
$ws.Range("N132").Value = -14740.4705
$ws.Range("H136").Value = 4472.923
$ws.Range("I136").Value = 3732.625
$ws.Range("K136").Value = 11197.875
$ws.Range("M136").Value = -8647.875
$ws.Range("H138").Value = 80000
$ws.Range("J138").Value = 80000
$ws.Range("L138").Value = 80000
$ws.Range("N138").Value = -90280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 55426.145
$ws.Range("J68").Value = 60180.5
$ws.Range("L68").Value = 60180.5
$ws.Range("N68").Value = -61802.5
$ws.Range("H71").Value = 55426.145
$ws.Range("J71").Value = 60180.5
$ws.Range("L71").Value = 180541.5
$ws.Range("N71").Value = -188653.5
$ws.Range("H122").Value = 1384.0546
$ws.Range("I122").Value = 1017.5714
$ws.Range("K122").Value = 3052.7142
$ws.Range("M122").Value = -602.7142000000003
$ws.Range("H132").Value = 913815.5600000001
$ws.Range("I132").Value = 1506782.2
$ws.Range("J132").Value = 4599.8667
$ws.Range("K132").Value = 4520346.6
$ws.Range("L132").Value = 13799.6001
$ws.Range("M132").Value = -4517816.6
$ws.Range("N132").Value = -18859.6001
$ws.Range("H136").Value = 11911278
$ws.Range("I136").Value = 15880338
$ws.Range("J136").Value = 4096.125
$ws.Range("K136").Value = 47641014
$ws.Range("L136").Value = 12288.375
$ws.Range("M136").Value = -47638464
$ws.Range("N136").Value = -17388.375
